# Swap the contents of columns C (codeforiati:group-code) and D
# (codeforiati:group-name), including the header row, on the single
# worksheet of the workbook. This reproduces the reordering seen in the
# shared-strings table where the "group-name" value for each group now
# precedes the "group-code" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

$rangeC = $ws.Range("C1:C$lastRow")
$rangeD = $ws.Range("D1:D$lastRow")

# Capture current values before overwriting either column. Value2 is used
# instead of Value because it reliably round-trips array data through the
# COM interop layer.
$valuesC = $rangeC.Value2
$valuesD = $rangeD.Value2

# Write the swapped values back.
$rangeC.Value2 = $valuesD
$rangeD.Value2 = $valuesC
